$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1968.3529
$ws.Range("I98").Value = 1596.9231
$ws.Range("J98").Value = 3175.5
$ws.Range("K98").Value = 1596.9231
$ws.Range("L98").Value = 3175.5
$ws.Range("M98").Value = -98.92309999999998
$ws.Range("N98").Value = -6171.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 503.5
$ws.Range("I107").Value = 465.30768
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 465.30768
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1454.69232
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1968.3529
$ws.Range("I122").Value = 1596.9231
$ws.Range("J122").Value = 3175.5
$ws.Range("K122").Value = 4790.7693
$ws.Range("L122").Value = 9526.5
$ws.Range("M122").Value = -2340.7693
$ws.Range("N122").Value = -14426.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 10126.346
$ws.Range("I132").Value = 9186.5625
$ws.Range("J132").Value = 11630
$ws.Range("K132").Value = 27559.6875
$ws.Range("L132").Value = 34890
$ws.Range("M132").Value = -25029.6875
$ws.Range("N132").Value = -39950

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 20835648
$ws.Range("I137").Value = 27780418
$ws.Range("J137").Value = 1333
$ws.Range("K137").Value = 83341254
$ws.Range("L137").Value = 3999
$ws.Range("M137").Value = -83338704
$ws.Range("N137").Value = -9099

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1960.5741
$ws.Range("I138").Value = 931.5769
$ws.Range("J138").Value = 2916.0715
$ws.Range("K138").Value = 2794.7307
$ws.Range("L138").Value = 8748.2145
$ws.Range("M138").Value = 2345.2693
$ws.Range("N138").Value = -19028.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7301.846
$ws.Range("I32").Value = 7589.298
$ws.Range("J32").Value = 4599.8
$ws.Range("K32").Value = 7589.298
$ws.Range("L32").Value = 4599.8
$ws.Range("M32").Value = -7302.298
$ws.Range("N32").Value = -5173.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3271.5833
$ws.Range("I45").Value = 2142.2144
$ws.Range("J45").Value = 4852.7
$ws.Range("K45").Value = 2142.2144
$ws.Range("L45").Value = 4852.7
$ws.Range("M45").Value = -1765.2144
$ws.Range("N45").Value = -5606.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2232.75
$ws.Range("I88").Value = 2064.3333
$ws.Range("J88").Value = 2401.1667
$ws.Range("K88").Value = 2064.3333
$ws.Range("L88").Value = 2401.1667
$ws.Range("M88").Value = -1658.3333
$ws.Range("N88").Value = -3213.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2232.75
$ws.Range("I91").Value = 2064.3333
$ws.Range("J91").Value = 2401.1667
$ws.Range("K91").Value = 2064.3333
$ws.Range("L91").Value = 2401.1667
$ws.Range("M91").Value = -660.3332999999998
$ws.Range("N91").Value = -5209.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2901.2307
$ws.Range("I122").Value = 1501.4546
$ws.Range("J122").Value = 10600
$ws.Range("K122").Value = 4504.3638
$ws.Range("L122").Value = 31800
$ws.Range("M122").Value = -2054.3638
$ws.Range("N122").Value = -36700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8774829
$ws.Range("I31").Value = 1759.4667
$ws.Range("J31").Value = 41673840
$ws.Range("K31").Value = 1759.4667
$ws.Range("L31").Value = 41673840
$ws.Range("M31").Value = -1464.4667
$ws.Range("N31").Value = -41674430

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8774829
$ws.Range("I34").Value = 1759.4667
$ws.Range("J34").Value = 41673840
$ws.Range("K34").Value = 1759.4667
$ws.Range("L34").Value = 41673840
$ws.Range("M34").Value = -1557.4667
$ws.Range("N34").Value = -41674244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1976.7368
$ws.Range("I94").Value = 3421.2
$ws.Range("J94").Value = 1460.8572
$ws.Range("K94").Value = 3421.2
$ws.Range("L94").Value = 1460.8572
$ws.Range("M94").Value = -2970.2
$ws.Range("N94").Value = -2362.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63.5
$ws.Range("I2").Value = 39.5
$ws.Range("J2").Value = 73.09999999999999
$ws.Range("K2").Value = 237
$ws.Range("L2").Value = 438.6
$ws.Range("M2").Value = -124
$ws.Range("N2").Value = -664.5999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1185.875
$ws.Range("I5").Value = 521.4
$ws.Range("J5").Value = 2293.3333
$ws.Range("K5").Value = 1564.2
$ws.Range("L5").Value = 6879.999899999999
$ws.Range("M5").Value = -1452.2
$ws.Range("N5").Value = -7103.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2450
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2611.111
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 7833.333
$ws.Range("M17").Value = -2831
$ws.Range("N17").Value = -8171.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 20430.908
$ws.Range("J34").Value = 24810
$ws.Range("L34").Value = 74430
$ws.Range("N34").Value = -74598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3628.7556
$ws.Range("I122").Value = 414.83334
$ws.Range("J122").Value = 3858.3215
$ws.Range("K122").Value = 3733.50006
$ws.Range("L122").Value = 34724.8935
$ws.Range("M122").Value = -1283.50006
$ws.Range("N122").Value = -39624.8935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3771.6667
$ws.Range("I125").Value = 2243.3333
$ws.Range("J125").Value = 5300
$ws.Range("K125").Value = 6729.999899999999
$ws.Range("L125").Value = 15900
$ws.Range("M125").Value = -1809.999899999999
$ws.Range("N125").Value = -25740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1633
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1949.5
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 17545.5
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -22605.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3943.9644
$ws.Range("I134").Value = 1500
$ws.Range("J134").Value = 4131.9614
$ws.Range("K134").Value = 4500
$ws.Range("L134").Value = 12395.8842
$ws.Range("M134").Value = 570
$ws.Range("N134").Value = -22535.8842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1185.875
$ws.Range("I135").Value = 521.4
$ws.Range("J135").Value = 2293.3333
$ws.Range("K135").Value = 4692.599999999999
$ws.Range("L135").Value = 20639.9997
$ws.Range("M135").Value = -2157.599999999999
$ws.Range("N135").Value = -25709.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2249.3845
$ws.Range("I132").Value = 2033.7142
$ws.Range("J132").Value = 2501
$ws.Range("K132").Value = 6101.142599999999
$ws.Range("L132").Value = 7503
$ws.Range("M132").Value = -3571.142599999999
$ws.Range("N132").Value = -12563

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2327.8462
$ws.Range("I82").Value = 1980.3334
$ws.Range("J82").Value = 2625.7144
$ws.Range("K82").Value = 1980.3334
$ws.Range("L82").Value = 2625.7144
$ws.Range("M82").Value = -1619.3334
$ws.Range("N82").Value = -3347.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2327.8462
$ws.Range("I85").Value = 1980.3334
$ws.Range("J85").Value = 2625.7144
$ws.Range("K85").Value = 1980.3334
$ws.Range("L85").Value = 2625.7144
$ws.Range("M85").Value = -732.3334
$ws.Range("N85").Value = -5121.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 37604.367
$ws.Range("I132").Value = 44269.12
$ws.Range("K132").Value = 132807.36
$ws.Range("M132").Value = -130277.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2413.4285
$ws.Range("I132").Value = 2031.3334
$ws.Range("J132").Value = 2795.524
$ws.Range("K132").Value = 6094.0002
$ws.Range("L132").Value = 8386.572
$ws.Range("M132").Value = -3564.0002
$ws.Range("N132").Value = -13446.572
